# Auto-generated PowerShell Word COM-interop script
$d = $word.ActiveDocument

# Update the title/date paragraph
$d.Content.Find.Execute("2024-02-28 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-29 Thursday", 2) | Out-Null

# Update table cell values (row-major order matches source document)
$t = $d.Tables.Item(1)

$newValues = @(
    @("70+21=91", "16+16=32", "60+32=92", "69-36=33", "11+58=69"),
    @("31-7=24", "4+45=49", "20+53=73", "67-58=9", "44+51=95"),
    @("86-79=7", "52-23=29", "70-54=16", "75-21=54", "78-11=67"),
    @("46-8=38", "9+62=71", "65-37=28", "10+79=89", "92-6=86"),
    @("47+23=70", "16+35=51", "68-58=10", "94-41=53", "78-77=1"),
    @("35+25=60", "46-38=8", "20+1=21", "15-4=11", "54-10=44"),
    @("89+1=90", "79+6=85", "28+63=91", "90-15=75", "14-12=2"),
    @("58+6=64", "37+32=69", "6+34=40", "64-43=21", "42-23=19"),
    @("55-30=25", "88-35=53", "5+3=8", "41+53=94", "45+10=55"),
    @("12+46=58", "77-47=30", "96-61=35", "46-44=2", "76-1=75"),
    @("42-3=39", "56+39=95", "54-26=28", "26+42=68", "3+60=63"),
    @("28-16=12", "90-37=53", "64-11=53", "52+40=92", "95-68=27"),
    @("69+14=83", "58-3=55", "79+9=88", "63-38=25", "91-2=89"),
    @("65-23=42", "74-42=32", "70+14=84", "72-65=7", "38+0=38"),
    @("78+21=99", "28-23=5", "18+35=53", "85-63=22", "19+76=95"),
    @("53-24=29", "77-1=76", "2+32=34", "62-53=9", "48-3=45"),
    @("38-29=9", "82-57=25", "68-29=39", "31+37=68", "36-4=32"),
    @("80-9=71", "33-30=3", "27+32=59", "29+44=73", "52+46=98"),
    @("39-21=18", "15+20=35", "16+82=98", "60-6=54", "57+25=82"),
    @("24-1=23", "95-50=45", "53-9=44", "65+11=76", "59+2=61")
)

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$r - 1][$c - 1]
    }
}

Write-Host "Done updating title and" ($t.Rows.Count * $t.Columns.Count) "table cells."
